$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, centered, bordered) from L1 onto the new header cells
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)

$ws.Cells.Item(1, 13).Value = "31/12/2023"
$ws.Cells.Item(1, 14).Value = "31/03/2024"
$ws.Cells.Item(1, 15).Value = "30/06/2024"
$ws.Cells.Item(2, 13).Value = 309631.008
$ws.Cells.Item(2, 14).Value = 300476.992
$ws.Cells.Item(2, 15).Value = 275648
$ws.Cells.Item(3, 13).Value = 113076
$ws.Cells.Item(3, 14).Value = 102870
$ws.Cells.Item(3, 15).Value = 96900
$ws.Cells.Item(4, 13).Value = 83457
$ws.Cells.Item(4, 14).Value = 74257
$ws.Cells.Item(4, 15).Value = 60438
$ws.Cells.Item(5, 13).Value = 1665
$ws.Cells.Item(5, 14).Value = 1700
$ws.Cells.Item(5, 15).Value = 1144
$ws.Cells.Item(6, 13).Value = 11853
$ws.Cells.Item(6, 14).Value = 11275
$ws.Cells.Item(6, 15).Value = 10834
$ws.Cells.Item(7, 13).Value = 0
$ws.Cells.Item(7, 14).Value = 0
$ws.Cells.Item(7, 15).Value = 0
$ws.Cells.Item(8, 13).Value = 0
$ws.Cells.Item(8, 14).Value = 0
$ws.Cells.Item(8, 15).Value = 0
$ws.Cells.Item(9, 13).Value = 4967
$ws.Cells.Item(9, 14).Value = 4810
$ws.Cells.Item(9, 15).Value = 5713
$ws.Cells.Item(10, 13).Value = 3959
$ws.Cells.Item(10, 14).Value = 4041
$ws.Cells.Item(10, 15).Value = 4430
$ws.Cells.Item(11, 13).Value = 7175
$ws.Cells.Item(11, 14).Value = 6787
$ws.Cells.Item(11, 15).Value = 14341
$ws.Cells.Item(12, 13).Value = 19175
$ws.Cells.Item(12, 14).Value = 18087
$ws.Cells.Item(12, 15).Value = 24706
$ws.Cells.Item(13, 13).Value = 0
$ws.Cells.Item(13, 14).Value = 0
$ws.Cells.Item(13, 15).Value = 0
$ws.Cells.Item(14, 13).Value = 0
$ws.Cells.Item(14, 14).Value = 0
$ws.Cells.Item(14, 15).Value = 0
$ws.Cells.Item(15, 13).Value = 0
$ws.Cells.Item(15, 14).Value = 0
$ws.Cells.Item(15, 15).Value = 0
$ws.Cells.Item(16, 13).Value = 0
$ws.Cells.Item(16, 14).Value = 0
$ws.Cells.Item(16, 15).Value = 0
$ws.Cells.Item(17, 13).Value = 0
$ws.Cells.Item(17, 14).Value = 0
$ws.Cells.Item(17, 15).Value = 0
$ws.Cells.Item(18, 13).Value = 0
$ws.Cells.Item(18, 14).Value = 0
$ws.Cells.Item(18, 15).Value = 0
$ws.Cells.Item(19, 13).Value = 11102
$ws.Cells.Item(19, 14).Value = 11193
$ws.Cells.Item(19, 15).Value = 12171
$ws.Cells.Item(20, 13).Value = 0
$ws.Cells.Item(20, 14).Value = 0
$ws.Cells.Item(20, 15).Value = 0
$ws.Cells.Item(21, 13).Value = 188
$ws.Cells.Item(21, 14).Value = 185
$ws.Cells.Item(21, 15).Value = 7060
$ws.Cells.Item(22, 13).Value = 76623
$ws.Cells.Item(22, 14).Value = 76327
$ws.Cells.Item(22, 15).Value = 60082
$ws.Cells.Item(23, 13).Value = 9813
$ws.Cells.Item(23, 14).Value = 11553
$ws.Cells.Item(23, 15).Value = 11086
$ws.Cells.Item(24, 13).Value = 90944
$ws.Cells.Item(24, 14).Value = 91640
$ws.Cells.Item(24, 15).Value = 82874
$ws.Cells.Item(25, 13).Value = 0
$ws.Cells.Item(25, 14).Value = 0
$ws.Cells.Item(25, 15).Value = 0
$ws.Cells.Item(26, 13).Value = 309631.008
$ws.Cells.Item(26, 14).Value = 300476.992
$ws.Cells.Item(26, 15).Value = 275648
$ws.Cells.Item(27, 13).Value = 23384
$ws.Cells.Item(27, 14).Value = 23615
$ws.Cells.Item(27, 15).Value = 21608
$ws.Cells.Item(28, 13).Value = 1933
$ws.Cells.Item(28, 14).Value = 1473
$ws.Cells.Item(28, 15).Value = 1309
$ws.Cells.Item(29, 13).Value = 0
$ws.Cells.Item(29, 14).Value = 0
$ws.Cells.Item(29, 15).Value = 0
$ws.Cells.Item(30, 13).Value = 843
$ws.Cells.Item(30, 14).Value = 1168
$ws.Cells.Item(30, 15).Value = 1696
$ws.Cells.Item(31, 13).Value = 547
$ws.Cells.Item(31, 14).Value = 694
$ws.Cells.Item(31, 15).Value = 539
$ws.Cells.Item(32, 13).Value = 0
$ws.Cells.Item(32, 14).Value = 0
$ws.Cells.Item(32, 15).Value = 0
$ws.Cells.Item(33, 13).Value = 0
$ws.Cells.Item(33, 14).Value = 0
$ws.Cells.Item(33, 15).Value = 0
$ws.Cells.Item(34, 13).Value = 20061
$ws.Cells.Item(34, 14).Value = 20280
$ws.Cells.Item(34, 15).Value = 18064
$ws.Cells.Item(35, 13).Value = 0
$ws.Cells.Item(35, 14).Value = 0
$ws.Cells.Item(35, 15).Value = 0
$ws.Cells.Item(36, 13).Value = 0
$ws.Cells.Item(36, 14).Value = 0
$ws.Cells.Item(36, 15).Value = 0
$ws.Cells.Item(37, 13).Value = 2470
$ws.Cells.Item(37, 14).Value = 4566
$ws.Cells.Item(37, 15).Value = 6083
$ws.Cells.Item(38, 13).Value = 0
$ws.Cells.Item(38, 14).Value = 2383
$ws.Cells.Item(38, 15).Value = 2595
$ws.Cells.Item(39, 13).Value = 0
$ws.Cells.Item(39, 14).Value = 0
$ws.Cells.Item(39, 15).Value = 0
$ws.Cells.Item(40, 13).Value = 0
$ws.Cells.Item(40, 14).Value = 0
$ws.Cells.Item(40, 15).Value = 0
$ws.Cells.Item(41, 13).Value = 0
$ws.Cells.Item(41, 14).Value = 0
$ws.Cells.Item(41, 15).Value = 0
$ws.Cells.Item(42, 13).Value = 0
$ws.Cells.Item(42, 14).Value = 0
$ws.Cells.Item(42, 15).Value = 0
$ws.Cells.Item(43, 13).Value = 2470
$ws.Cells.Item(43, 14).Value = 2183
$ws.Cells.Item(43, 15).Value = 3488
$ws.Cells.Item(44, 13).Value = 0
$ws.Cells.Item(44, 14).Value = 0
$ws.Cells.Item(44, 15).Value = 0
$ws.Cells.Item(45, 13).Value = 0
$ws.Cells.Item(45, 14).Value = 0
$ws.Cells.Item(45, 15).Value = 0
$ws.Cells.Item(46, 13).Value = -172
$ws.Cells.Item(46, 14).Value = -191
$ws.Cells.Item(46, 15).Value = -12
$ws.Cells.Item(47, 13).Value = 283948.992
$ws.Cells.Item(47, 14).Value = 272487
$ws.Cells.Item(47, 15).Value = 247968.992
$ws.Cells.Item(48, 13).Value = 581164.032
$ws.Cells.Item(48, 14).Value = 581164.032
$ws.Cells.Item(48, 15).Value = 581164.032
$ws.Cells.Item(49, 13).Value = 20568
$ws.Cells.Item(49, 14).Value = 20568
$ws.Cells.Item(49, 15).Value = 20568
$ws.Cells.Item(50, 13).Value = 0
$ws.Cells.Item(50, 14).Value = 0
$ws.Cells.Item(50, 15).Value = 0
$ws.Cells.Item(51, 13).Value = 0
$ws.Cells.Item(51, 14).Value = 0
$ws.Cells.Item(51, 15).Value = 0
$ws.Cells.Item(52, 13).Value = -317783.008
$ws.Cells.Item(52, 14).Value = -329244.992
$ws.Cells.Item(52, 15).Value = -353763.008
$ws.Cells.Item(53, 13).Value = 0
$ws.Cells.Item(53, 14).Value = 0
$ws.Cells.Item(53, 15).Value = 0
$ws.Cells.Item(54, 13).Value = 0
$ws.Cells.Item(54, 14).Value = 0
$ws.Cells.Item(54, 15).Value = 0
$ws.Cells.Item(55, 13).Value = 0
$ws.Cells.Item(55, 14).Value = 0
$ws.Cells.Item(55, 15).Value = 0
$ws.Cells.Item(56, 13).Value = 0
$ws.Cells.Item(56, 14).Value = 0
$ws.Cells.Item(56, 15).Value = 0
$ws.Cells.Item(59, 13).Value = 13405
$ws.Cells.Item(59, 14).Value = 11386
$ws.Cells.Item(59, 15).Value = 13070
$ws.Cells.Item(60, 13).Value = -8972
$ws.Cells.Item(60, 14).Value = -9148
$ws.Cells.Item(60, 15).Value = -8609
$ws.Cells.Item(61, 13).Value = 4433
$ws.Cells.Item(61, 14).Value = 2238
$ws.Cells.Item(61, 15).Value = 4461
$ws.Cells.Item(62, 13).Value = -3130
$ws.Cells.Item(62, 14).Value = -2214
$ws.Cells.Item(62, 15).Value = -1776
$ws.Cells.Item(63, 13).Value = -18577
$ws.Cells.Item(63, 14).Value = -12410
$ws.Cells.Item(63, 15).Value = -16371
$ws.Cells.Item(64, 13).Value = 0
$ws.Cells.Item(64, 14).Value = 0
$ws.Cells.Item(64, 15).Value = 0
$ws.Cells.Item(65, 13).Value = -128884
$ws.Cells.Item(65, 14).Value = 0
$ws.Cells.Item(65, 15).Value = -12082
$ws.Cells.Item(66, 13).Value = -7048
$ws.Cells.Item(66, 14).Value = 27
$ws.Cells.Item(66, 15).Value = 0
$ws.Cells.Item(67, 13).Value = -262
$ws.Cells.Item(67, 14).Value = -40
$ws.Cells.Item(67, 15).Value = 1654
$ws.Cells.Item(68, 13).Value = -185
$ws.Cells.Item(68, 14).Value = 1475
$ws.Cells.Item(68, 15).Value = 488
$ws.Cells.Item(69, 13).Value = 2727
$ws.Cells.Item(69, 14).Value = 2476
$ws.Cells.Item(69, 15).Value = 2250
$ws.Cells.Item(70, 13).Value = -2912
$ws.Cells.Item(70, 14).Value = -1001
$ws.Cells.Item(70, 15).Value = -1762
$ws.Cells.Item(74, 13).Value = -153652.992
$ws.Cells.Item(74, 14).Value = -10924
$ws.Cells.Item(74, 15).Value = -23626
$ws.Cells.Item(75, 13).Value = -766
$ws.Cells.Item(75, 14).Value = -665
$ws.Cells.Item(75, 15).Value = -678
$ws.Cells.Item(76, 13).Value = -35204
$ws.Cells.Item(76, 14).Value = 43
$ws.Cells.Item(76, 15).Value = 279
$ws.Cells.Item(79, 13).Value = 215
$ws.Cells.Item(79, 14).Value = 84
$ws.Cells.Item(79, 15).Value = -245
$ws.Cells.Item(80, 13).Value = -189408
$ws.Cells.Item(80, 14).Value = -11462
$ws.Cells.Item(80, 15).Value = -24270
